# Apply the commit "Added a few more slots":
#  1. Insert a new "Meta description" paragraph right after the title
#     (Heading1) paragraph.
#  2. At the end of the document, delete the duplicated bold title
#     paragraph and replace the italic paragraph's text with the new
#     AI image-generation prompt.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the "Meta description" paragraph after the title paragraph.
# ---------------------------------------------------------------------
$titlePar = $d.Paragraphs.Item(1)
$titlePar.Range.InsertParagraphAfter()

$metaPar = $d.Paragraphs.Item(2)
$metaPar.Style = "Normal"

$metaRange = $metaPar.Range
$metaRange.Collapse(1)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaRange.InsertXML($metaXml) | Out-Null

$metaPar = $d.Paragraphs.Item(2)
$insertAt = $metaPar.Range.End - 1
$tailRange = $d.Range($insertAt, $insertAt)
$tailRange.InsertAfter(": Discover the pros and cons of Caishen" + [char]0x2019 + "s Fortune XL video slot game. Play for free and enjoy medium volatility and a 96% RTP.")

# ---------------------------------------------------------------------
# 2) At the end of the document: drop the duplicated bold title
#    paragraph, and rewrite the italic paragraph's text.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPar = $d.Paragraphs.Item($count)
$boldPar = $d.Paragraphs.Item($count - 1)

$boldPar.Range.Delete()

$lastPar = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPar.Range.Text = "Create a feature image for " + [char]0x22 + "Caishen's Fortune XL" + [char]0x22 + " that captures the vibrant and luxurious Chinese theme of the game. The image should be in a cartoon style and feature a happy Maya warrior wearing glasses, symbolizing the exciting adventure and good fortune that players can experience in this online slot game. The background of the image should showcase traditional Chinese designs and colors, with cherry blossoms and gazebos. The image should be visually stunning and eye-catching, drawing players to try their luck with Cai Shen's Fortune XL."
